$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.374.11"
$ws.Range("E2").Value = "  -0.16%  "

# Row 3
$ws.Range("D3").Value = "2.182.04"
$ws.Range("E3").Value = "  -1.36%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").Value = "'253.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.60%  "

# Row 6
$ws.Range("E6").Value = "  -0.56%  "

# Row 7
$ws.Range("D7").Value = "'73.29"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.34%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("E9").Value = "  -2.35%  "

# Row 10
$ws.Range("D10").Value = "'40.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.87%  "

# Row 11
$ws.Range("D11").Value = "'0.0914"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.86%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.101"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.07%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'6.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.20%  "

# Row 14
$ws.Range("D14").Value = "2.511.36"
$ws.Range("E14").Value = "  -1.26%  "

# Row 15
$ws.Range("D15").Value = "'14.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.39%  "

# Row 16
$ws.Range("D16").Value = "2.190.94"
$ws.Range("E16").Value = "  -0.79%  "

# Row 17
$ws.Range("D17").Value = "'0.770"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.58%  "

# Row 18
$ws.Range("D18").Value = "42.322.34"
$ws.Range("E18").Value = "  +0.05%  "

# Row 19
$ws.Range("E19").Value = "  -2.65%  "

# Row 20
$ws.Range("D20").Value = "'70.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.25%  "

# Row 21
$ws.Range("E21").Value = "  -0.63%  "

# Row 22
$ws.Range("D22").Value = "'226.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.78%  "

# Row 23
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "'9.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.33%  "

# Row 24
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").Value = "'2.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.92%  "

# Row 26
$ws.Range("D26").Value = "'10.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.16%  "

# Row 27
$ws.Range("D27").Value = "'3.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.13%  "

# Row 28
$ws.Range("E28").Value = "  +1.61%  "

# Row 29
$ws.Range("E29").Value = "  -1.99%  "

# Row 30
$ws.Range("D30").Value = "'170.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.22%  "

# Row 31
$ws.Range("D31").Value = "'36.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.23%  "

# Row 32
$ws.Range("D32").Value = "'20.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.94%  "

# Row 33
$ws.Range("D33").Value = "'0.0805"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.59%  "

# Row 34
$ws.Range("E34").Value = "  -4.30%  "

# Row 35
$ws.Range("E35").Value = "  -1.24%  "

# Row 36
$ws.Range("E36").Value = "  +0.12%  "

# Row 37
$ws.Range("D37").Value = "'4.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.48%  "

# Row 38
$ws.Range("D38").Value = "'0.0338"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.69%  "

# Row 39
$ws.Range("D39").Value = "'11.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.20%  "

# Row 40
$ws.Range("E40").Value = "  -3.36%  "

# Row 41
$ws.Range("E41").Value = "  +0.12%  "

# Row 42
$ws.Range("D42").Value = "'59.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.77%  "

# Row 43
$ws.Range("E43").Value = "  -6.20%  "

# Row 44
$ws.Range("D44").Value = "'102.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.79%  "

# Row 45
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'2.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.43%  "

# Row 46
$ws.Range("B46").Value = "WOONetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D46").Value = "'0.466"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.22%  "

# Row 47
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'8.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.40%  "

# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.0969"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.83%  "

# Row 49
$ws.Range("E49").Value = "  -1.00%  "

# Row 50
$ws.Range("D50").Value = "'1.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.93%  "

# Row 51
$ws.Range("E51").Value = "  +0.45%  "
